# Update cryptos list on Mon Mar 13 17:14:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to text format temporarily so numeric-looking
# strings (e.g. "1.003", "308.05") are not auto-converted to numbers by
# Excel; format is restored to General afterwards to match the source file.
# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.961.86"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +16.52%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.661.33"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +12.53%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.73%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.05"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +11.26%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9969"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.93%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3707"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +4.79%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3437"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +11.85%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.78"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +21.37%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.165"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +7.79%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07221"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +8.65%  "

# Row 12 - BinanceUSD
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9991"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.54%  "

# Row 13 - Solana
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.49"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +13.05%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.016"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +10.17%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.721"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +8.88%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.661.45"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +12.55%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +7.87%  "

# Row 18 - Dai
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9965"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +3.92%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06709"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +12.06%  "

# Row 20 - Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.45"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +17.72%  "

# Row 21 - Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.37"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +12.93%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.115"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +11.54%  "

# Row 23 - Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.96"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +7.97%  "

# Row 24 - WrappedBTC
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.920.81"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +16.09%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.388"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.22%  "

# Rows 26-27 swapped: LidoDAOToken <-> LEO (with updated price/volume)
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.379"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -8.13%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.663"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +27.87%  "

# Row 28 - Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.56"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.90%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.48"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +13.40%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.845.89"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +12.81%  "

# Row 31 - BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.78"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +10.58%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.350"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +28.61%  "

# Row 33 - HuobiToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.104"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +4.94%  "

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9744"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +22.27%  "

# Row 35 - WEMIXTOKEN
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.740"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +21.34%  "

# Row 36 - Stellar
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08381"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +5.90%  "

# Rows 37-38 swapped: FraxShare <-> Aptos (with updated price/volume)
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.25"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +19.83%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.944"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +21.94%  "

# Rows 39-40 swapped: InternetComputer(DFINITY) <-> Hedera (with updated price/volume)
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06341"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +11.75%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.296"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +12.69%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.280"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +7.45%  "

# Row 42 - VeChain
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02310"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +14.68%  "

# Row 43 - Algorand
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2075"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +12.55%  "

# Row 44 - TheSandbox
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6082"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +16.65%  "

# Row 45 - Frax
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9955"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +3.74%  "

# Row 46 - PancakeSwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.820"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +8.83%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.25"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +9.86%  "

# Row 48 - Decentraland
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5932"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +15.03%  "

# Row 49 - Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.87"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +6.08%  "

# Row 50 - NEARProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.994"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +10.37%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07074"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +10.50%  "
